$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.769.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.697.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.56%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.402"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "30.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.54%  "
$ws.Range("E14").Value = "  +8.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.179.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.622.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.697.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "358.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000107"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.13%  "
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("E27").Value = "  +3.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.68%  "
$ws.Range("E29").Value = "  +4.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "547.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.76%  "
$ws.Range("E31").Value = "  +4.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.47%  "
$ws.Range("E34").Value = "  +6.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0619"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.61%  "
$ws.Range("E47").Value = "  +2.66%  "
$ws.Range("E48").Value = "  +4.57%  "
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0994"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.47%  "
